$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.032993782094195
$ws.Cells.Item(2, 4).Value = 1.042341109638337
$ws.Cells.Item(2, 5).Value = 1.050632533155241
$ws.Cells.Item(2, 6).Value = 1.055486803584468
$ws.Cells.Item(2, 9).Value = 1.0395884734699
$ws.Cells.Item(2, 10).Value = 1.038120948107042
$ws.Cells.Item(2, 11).Value = 1.045118071181544
$ws.Cells.Item(2, 12).Value = 1.053386270475264
$ws.Cells.Item(2, 13).Value = 1.058227127627049
$ws.Cells.Item(2, 14).Value = 1.016697893379302

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033748985781565
$ws.Cells.Item(3, 4).Value = 1.042942308748059
$ws.Cells.Item(3, 5).Value = 1.051480690197687
$ws.Cells.Item(3, 6).Value = 1.056295623587447
$ws.Cells.Item(3, 9).Value = 1.039763583382552
$ws.Cells.Item(3, 10).Value = 1.0385196655569
$ws.Cells.Item(3, 11).Value = 1.045530631591462
$ws.Cells.Item(3, 12).Value = 1.054046820891078
$ws.Cells.Item(3, 13).Value = 1.058849410388343
$ws.Cells.Item(3, 14).Value = 1.016829985723613

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034238149969269
$ws.Cells.Item(4, 4).Value = 1.043331742380676
$ws.Cells.Item(4, 5).Value = 1.052030900535779
$ws.Cells.Item(4, 6).Value = 1.056820052591245
$ws.Cells.Item(4, 9).Value = 1.039875764098196
$ws.Cells.Item(4, 10).Value = 1.038777463300584
$ws.Cells.Item(4, 11).Value = 1.045797306163493
$ws.Cells.Item(4, 12).Value = 1.054474989457717
$ws.Cells.Item(4, 13).Value = 1.059252481179285
$ws.Cells.Item(4, 14).Value = 1.01691538100809

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034443911693111
$ws.Cells.Item(5, 4).Value = 1.043495558358169
$ws.Cells.Item(5, 5).Value = 1.052262541109882
$ws.Cells.Item(5, 6).Value = 1.057040776218326
$ws.Cells.Item(5, 9).Value = 1.039922654180576
$ws.Cells.Item(5, 10).Value = 1.038885792529208
$ws.Cells.Item(5, 11).Value = 1.045909347833794
$ws.Cells.Item(5, 12).Value = 1.054655168936469
$ws.Cells.Item(5, 13).Value = 1.059422029083664
$ws.Cells.Item(5, 14).Value = 1.016951262234843

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.034478466804625
$ws.Cells.Item(6, 4).Value = 1.043523069477728
$ws.Cells.Item(6, 5).Value = 1.052301453998919
$ws.Cells.Item(6, 6).Value = 1.057077851501084
$ws.Cells.Item(6, 9).Value = 1.039930511333359
$ws.Cells.Item(6, 10).Value = 1.038903978576214
$ws.Cells.Item(6, 11).Value = 1.045928156076707
$ws.Cells.Item(6, 12).Value = 1.054685432217195
$ws.Cells.Item(6, 13).Value = 1.059450502556919
$ws.Cells.Item(6, 14).Value = 1.016957285726427

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034240898909705
$ws.Cells.Item(7, 4).Value = 1.043333930914795
$ws.Cells.Item(7, 5).Value = 1.052033994425942
$ws.Cells.Item(7, 6).Value = 1.056823000917025
$ws.Cells.Item(7, 9).Value = 1.039876391710473
$ws.Cells.Item(7, 10).Value = 1.038778910994953
$ws.Cells.Item(7, 11).Value = 1.045798803539963
$ws.Cells.Item(7, 12).Value = 1.054477396329784
$ws.Cells.Item(7, 13).Value = 1.059254746307951
$ws.Cells.Item(7, 14).Value = 1.016915860529823

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033248902710588
$ws.Cells.Item(8, 4).Value = 1.042544200083843
$ws.Cells.Item(8, 5).Value = 1.050918881504929
$ws.Cells.Item(8, 6).Value = 1.055759925615657
$ws.Cells.Item(8, 9).Value = 1.039647885543254
$ws.Cells.Item(8, 10).Value = 1.038255736995736
$ws.Cells.Item(8, 11).Value = 1.045257554965678
$ws.Cells.Item(8, 12).Value = 1.053609350607476
$ws.Cells.Item(8, 13).Value = 1.0584373445947
$ws.Cells.Item(8, 14).Value = 1.016742550341839

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031504766637116
$ws.Cells.Item(9, 4).Value = 1.041155875669564
$ws.Cells.Item(9, 5).Value = 1.048964693100878
$ws.Cells.Item(9, 6).Value = 1.053894926396237
$ws.Cells.Item(9, 9).Value = 1.039236634179115
$ws.Cells.Item(9, 10).Value = 1.037332369034803
$ws.Cells.Item(9, 11).Value = 1.04430172554351
$ws.Cells.Item(9, 12).Value = 1.052085550412489
$ws.Cells.Item(9, 13).Value = 1.057000212764267
$ws.Cells.Item(9, 14).Value = 1.016436583375264

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.030344734323767
$ws.Cells.Item(10, 4).Value = 1.040232646123401
$ws.Cells.Item(10, 5).Value = 1.047669274046723
$ws.Cells.Item(10, 6).Value = 1.052657282823098
$ws.Cells.Item(10, 9).Value = 1.038956739548536
$ws.Cells.Item(10, 10).Value = 1.036715879169955
$ws.Cells.Item(10, 11).Value = 1.043663189804051
$ws.Cells.Item(10, 12).Value = 1.051073691306056
$ws.Cells.Item(10, 13).Value = 1.056044403841938
$ws.Cells.Item(10, 14).Value = 1.016232247250832

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.029843097314044
$ws.Cells.Item(11, 4).Value = 1.039833451814109
$ws.Cells.Item(11, 5).Value = 1.047110118251018
$ws.Cells.Item(11, 6).Value = 1.05212274546193
$ws.Cells.Item(11, 9).Value = 1.038834194514956
$ws.Cells.Item(11, 10).Value = 1.036448732208414
$ws.Cells.Item(11, 11).Value = 1.043386402060207
$ws.Cells.Item(11, 12).Value = 1.050636517568924
$ws.Cells.Item(11, 13).Value = 1.055631090645125
$ws.Cells.Item(11, 14).Value = 1.016143688130513

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.029656868739897
$ws.Cells.Item(12, 4).Value = 1.0396852608849
$ws.Cells.Item(12, 5).Value = 1.046902690825448
$ws.Cells.Item(12, 6).Value = 1.051924402545038
$ws.Cells.Item(12, 9).Value = 1.038788474046301
$ws.Cells.Item(12, 10).Value = 1.036349472762644
$ws.Cells.Item(12, 11).Value = 1.043283547470206
$ws.Cells.Item(12, 12).Value = 1.050474278975288
$ws.Cells.Item(12, 13).Value = 1.05547765340027
$ws.Cells.Item(12, 14).Value = 1.016110781730031

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02969681076658
$ws.Cells.Item(13, 4).Value = 1.039717044346413
$ws.Cells.Item(13, 5).Value = 1.046947172553611
$ws.Cells.Item(13, 6).Value = 1.051966938324195
$ws.Cells.Item(13, 9).Value = 1.038798290359546
$ws.Cells.Item(13, 10).Value = 1.03637076554989
$ws.Cells.Item(13, 11).Value = 1.043305612065669
$ws.Cells.Item(13, 12).Value = 1.050509073014841
$ws.Cells.Item(13, 13).Value = 1.055510562304255
$ws.Cells.Item(13, 14).Value = 1.016117840782291

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.029827701515463
$ws.Cells.Item(14, 4).Value = 1.03982120050717
$ws.Cells.Item(14, 5).Value = 1.047092966757234
$ws.Cells.Item(14, 6).Value = 1.052106346117705
$ws.Cells.Item(14, 9).Value = 1.038830419358187
$ws.Cells.Item(14, 10).Value = 1.036440527979219
$ws.Cells.Item(14, 11).Value = 1.043377900942848
$ws.Cells.Item(14, 12).Value = 1.050623103863107
$ws.Cells.Item(14, 13).Value = 1.055618405716422
$ws.Cells.Item(14, 14).Value = 1.016140968311202

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.029908361164825
$ws.Cells.Item(15, 4).Value = 1.039885386222982
$ws.Cells.Item(15, 5).Value = 1.047182830945391
$ws.Cells.Item(15, 6).Value = 1.052192267494825
$ws.Cells.Item(15, 9).Value = 1.038850188381565
$ws.Cells.Item(15, 10).Value = 1.036483507087885
$ws.Cells.Item(15, 11).Value = 1.043422434810555
$ws.Cells.Item(15, 12).Value = 1.050693381590441
$ws.Cells.Item(15, 13).Value = 1.055684863008789
$ws.Cells.Item(15, 14).Value = 1.016155216421497

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.030378040484665
$ws.Cells.Item(16, 4).Value = 1.040259151514397
$ws.Cells.Item(16, 5).Value = 1.047706420832307
$ws.Cells.Item(16, 6).Value = 1.05269278736611
$ws.Cells.Item(16, 9).Value = 1.038964844114034
$ws.Cells.Item(16, 10).Value = 1.036733604672248
$ws.Cells.Item(16, 11).Value = 1.043681553145962
$ws.Cells.Item(16, 12).Value = 1.051102725632083
$ws.Cells.Item(16, 13).Value = 1.056071845994694
$ws.Cells.Item(16, 14).Value = 1.01623812297636

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03067283719122
$ws.Cells.Item(17, 4).Value = 1.04049375863636
$ws.Cells.Item(17, 5).Value = 1.048035330105477
$ws.Cells.Item(17, 6).Value = 1.053007118698505
$ws.Cells.Item(17, 9).Value = 1.039036404115319
$ws.Cells.Item(17, 10).Value = 1.036890430904303
$ws.Cells.Item(17, 11).Value = 1.043844012561205
$ws.Cells.Item(17, 12).Value = 1.051359756732107
$ws.Cells.Item(17, 13).Value = 1.056314740954145
$ws.Cells.Item(17, 14).Value = 1.016290106908387

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.030844851123488
$ws.Cells.Item(18, 4).Value = 1.040630655884656
$ws.Cells.Item(18, 5).Value = 1.048227347879415
$ws.Cells.Item(18, 6).Value = 1.053190594932554
$ws.Cells.Item(18, 9).Value = 1.039078013655116
$ws.Cells.Item(18, 10).Value = 1.036981885339803
$ws.Cells.Item(18, 11).Value = 1.043938743585544
$ws.Cells.Item(18, 12).Value = 1.051509771896536
$ws.Cells.Item(18, 13).Value = 1.056456471229539
$ws.Cells.Item(18, 14).Value = 1.016320420485024

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.030903514227246
$ws.Cells.Item(19, 4).Value = 1.040677343558396
$ws.Cells.Item(19, 5).Value = 1.048292849843183
$ws.Cells.Item(19, 6).Value = 1.053253177925208
$ws.Cells.Item(19, 9).Value = 1.039092179322962
$ws.Cells.Item(19, 10).Value = 1.037013065566082
$ws.Cells.Item(19, 11).Value = 1.043971039488316
$ws.Cells.Item(19, 12).Value = 1.051560938957433
$ws.Cells.Item(19, 13).Value = 1.056504806654126
$ws.Cells.Item(19, 14).Value = 1.016330755299309

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.030641201660973
$ws.Cells.Item(20, 4).Value = 1.040468581821579
$ws.Cells.Item(20, 5).Value = 1.048000023596092
$ws.Cells.Item(20, 6).Value = 1.052973380243489
$ws.Cells.Item(20, 9).Value = 1.039028739865256
$ws.Cells.Item(20, 10).Value = 1.03687360695718
$ws.Cells.Item(20, 11).Value = 1.043826585172637
$ws.Cells.Item(20, 12).Value = 1.051332170067181
$ws.Cells.Item(20, 13).Value = 1.056288675050691
$ws.Cells.Item(20, 14).Value = 1.016284530325013

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.029789154628377
$ws.Cells.Item(21, 4).Value = 1.039790526679717
$ws.Cells.Item(21, 5).Value = 1.047050026584713
$ws.Cells.Item(21, 6).Value = 1.05206528823636
$ws.Cells.Item(21, 9).Value = 1.038820963732302
$ws.Cells.Item(21, 10).Value = 1.036419985477383
$ws.Cells.Item(21, 11).Value = 1.04335661485422
$ws.Cells.Item(21, 12).Value = 1.050589520542802
$ws.Cells.Item(21, 13).Value = 1.055586646139931
$ws.Cells.Item(21, 14).Value = 1.016134158146957

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029254028453254
$ws.Cells.Item(22, 4).Value = 1.039364714354379
$ws.Cells.Item(22, 5).Value = 1.04645427699239
$ws.Cells.Item(22, 6).Value = 1.051495539585246
$ws.Cells.Item(22, 9).Value = 1.038689159579011
$ws.Cells.Item(22, 10).Value = 1.036134607523721
$ws.Cells.Item(22, 11).Value = 1.043060876048474
$ws.Cells.Item(22, 12).Value = 1.050123439829434
$ws.Cells.Item(22, 13).Value = 1.055145749299079
$ws.Cells.Item(22, 14).Value = 1.016039546299802

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.029537652443149
$ws.Cells.Item(23, 4).Value = 1.039590396736256
$ws.Cells.Item(23, 5).Value = 1.04676994737291
$ws.Cells.Item(23, 6).Value = 1.051797459199398
$ws.Cells.Item(23, 9).Value = 1.038759141791218
$ws.Cells.Item(23, 10).Value = 1.036285907280292
$ws.Cells.Item(23, 11).Value = 1.043217675995157
$ws.Cells.Item(23, 12).Value = 1.050370436632194
$ws.Cells.Item(23, 13).Value = 1.055379429387222
$ws.Cells.Item(23, 14).Value = 1.016089708019265

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.030655496197846
$ws.Cells.Item(24, 4).Value = 1.040479957971133
$ws.Cells.Item(24, 5).Value = 1.048015976561214
$ws.Cells.Item(24, 6).Value = 1.052988624791218
$ws.Cells.Item(24, 9).Value = 1.039032203412354
$ws.Cells.Item(24, 10).Value = 1.036881209035596
$ws.Cells.Item(24, 11).Value = 1.043834459948559
$ws.Cells.Item(24, 12).Value = 1.051344635005644
$ws.Cells.Item(24, 13).Value = 1.056300452944608
$ws.Cells.Item(24, 14).Value = 1.016287050167069

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.031955195436752
$ws.Cells.Item(25, 4).Value = 1.041514390162861
$ws.Cells.Item(25, 5).Value = 1.049468607092498
$ws.Cells.Item(25, 6).Value = 1.054376079808249
$ws.Cells.Item(25, 9).Value = 1.039343965703351
$ws.Cells.Item(25, 10).Value = 1.037571247965551
$ws.Cells.Item(25, 11).Value = 1.044549068374885
$ws.Cells.Item(25, 12).Value = 1.052478790460669
$ws.Cells.Item(25, 13).Value = 1.057371351742378
$ws.Cells.Item(25, 14).Value = 1.01651574824248
